$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.033524121305692
$ws.Cells.Item(2, 4).Value = 1.043711246246137
$ws.Cells.Item(2, 5).Value = 1.051408696855093
$ws.Cells.Item(2, 6).Value = 1.05650988762704
$ws.Cells.Item(2, 9).Value = 1.038235441160286
$ws.Cells.Item(2, 10).Value = 1.03864837572425
$ws.Cells.Item(2, 11).Value = 1.046484342436077
$ws.Cells.Item(2, 12).Value = 1.054160280470766
$ws.Cells.Item(2, 13).Value = 1.059247401728483
$ws.Cells.Item(2, 14).Value = 1.016859921101848
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.034286375530582
$ws.Cells.Item(3, 4).Value = 1.044330046317329
$ws.Cells.Item(3, 5).Value = 1.052267596531864
$ws.Cells.Item(3, 6).Value = 1.057331748497834
$ws.Cells.Item(3, 9).Value = 1.038392259226163
$ws.Cells.Item(3, 10).Value = 1.039054328792287
$ws.Cells.Item(3, 11).Value = 1.046914735619861
$ws.Cells.Item(3, 12).Value = 1.054831701528539
$ws.Cells.Item(3, 13).Value = 1.059882894866289
$ws.Cells.Item(3, 14).Value = 1.016995203377904
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.034780233276087
$ws.Cells.Item(4, 4).Value = 1.044731032906577
$ws.Cells.Item(4, 5).Value = 1.052824835392386
$ws.Cells.Item(4, 6).Value = 1.057864748094505
$ws.Cells.Item(4, 9).Value = 1.038492854421043
$ws.Cells.Item(4, 10).Value = 1.039316941390243
$ws.Cells.Item(4, 11).Value = 1.047193116403528
$ws.Cells.Item(4, 12).Value = 1.055266983402796
$ws.Cells.Item(4, 13).Value = 1.060294648198379
$ws.Cells.Item(4, 14).Value = 1.017082695233123
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.034987999628495
$ws.Cells.Item(5, 4).Value = 1.044899744990607
$ws.Cells.Item(5, 5).Value = 1.053059449267036
$ws.Cells.Item(5, 6).Value = 1.058089106466545
$ws.Cells.Item(5, 9).Value = 1.038534934052772
$ws.Cells.Item(5, 10).Value = 1.039427326612931
$ws.Cells.Item(5, 11).Value = 1.047310119384432
$ws.Cells.Item(5, 12).Value = 1.05545017228951
$ws.Cells.Item(5, 13).Value = 1.060467878451152
$ws.Cells.Item(5, 14).Value = 1.017119465622807
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.035022893170703
$ws.Cells.Item(6, 4).Value = 1.044928080489048
$ws.Cells.Item(6, 5).Value = 1.053098862471737
$ws.Cells.Item(6, 6).Value = 1.058126793896005
$ws.Cells.Item(6, 9).Value = 1.038541987037821
$ws.Cells.Item(6, 10).Value = 1.03944585973709
$ws.Cells.Item(6, 11).Value = 1.047329762999641
$ws.Cells.Item(6, 12).Value = 1.055480941975556
$ws.Cells.Item(6, 13).Value = 1.06049697209557
$ws.Cells.Item(6, 14).Value = 1.017125638863817
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.034783008878857
$ws.Cells.Item(7, 4).Value = 1.044733286708517
$ws.Cells.Item(7, 5).Value = 1.052827968939204
$ws.Cells.Item(7, 6).Value = 1.057867744863489
$ws.Cells.Item(7, 9).Value = 1.038493417519523
$ws.Cells.Item(7, 10).Value = 1.039318416431827
$ws.Cells.Item(7, 11).Value = 1.047194679916037
$ws.Cells.Item(7, 12).Value = 1.055269430412172
$ws.Cells.Item(7, 13).Value = 1.060296962404208
$ws.Cells.Item(7, 14).Value = 1.017083186605396
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033781597480278
$ws.Cells.Item(8, 4).Value = 1.043920251065131
$ws.Cells.Item(8, 5).Value = 1.05169865971801
$ws.Cells.Item(8, 6).Value = 1.056787389490717
$ws.Cells.Item(8, 9).Value = 1.038288619620223
$ws.Cells.Item(8, 10).Value = 1.038785582404612
$ws.Cells.Item(8, 11).Value = 1.04662981825904
$ws.Cells.Item(8, 12).Value = 1.054387017969378
$ws.Cells.Item(8, 13).Value = 1.059462055344869
$ws.Cells.Item(8, 14).Value = 1.016905649387706
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.032021876882297
$ws.Cells.Item(9, 4).Value = 1.042492116750668
$ws.Cells.Item(9, 5).Value = 1.049720047886395
$ws.Cells.Item(9, 6).Value = 1.054892953760097
$ws.Cells.Item(9, 9).Value = 1.037921058500404
$ws.Cells.Item(9, 10).Value = 1.037846210051273
$ws.Cells.Item(9, 11).Value = 1.045633660840718
$ws.Cells.Item(9, 12).Value = 1.052838509217487
$ws.Cells.Item(9, 13).Value = 1.05799510746232
$ws.Cells.Item(9, 14).Value = 1.016592482969039
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.03085213154388
$ws.Cells.Item(10, 4).Value = 1.041543192552042
$ws.Cells.Item(10, 5).Value = 1.048408744958113
$ws.Cells.Item(10, 6).Value = 1.053636367272013
$ws.Cells.Item(10, 9).Value = 1.037671569745549
$ws.Cells.Item(10, 10).Value = 1.037219734686705
$ws.Cells.Item(10, 11).Value = 1.044969102162975
$ws.Cells.Item(10, 12).Value = 1.051810586614146
$ws.Cells.Item(10, 13).Value = 1.057020114409214
$ws.Cells.Item(10, 14).Value = 1.016383515299106
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.030346449243923
$ws.Cells.Item(11, 4).Value = 1.0411330728517
$ws.Cells.Item(11, 5).Value = 1.047842805834544
$ws.Cells.Item(11, 6).Value = 1.05309378912461
$ws.Cells.Item(11, 9).Value = 1.037562492651705
$ws.Cells.Item(11, 10).Value = 1.036948426135298
$ws.Cells.Item(11, 11).Value = 1.044681250981125
$ws.Cells.Item(11, 12).Value = 1.051366554307569
$ws.Cells.Item(11, 13).Value = 1.056598658797985
$ws.Cells.Item(11, 14).Value = 1.016292990673146
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.030158742194022
$ws.Cells.Item(12, 4).Value = 1.040980853898019
$ws.Cells.Item(12, 5).Value = 1.047632872983634
$ws.Cells.Item(12, 6).Value = 1.05289248379388
$ws.Cells.Item(12, 9).Value = 1.037521820030239
$ws.Cells.Item(12, 10).Value = 1.036847645260756
$ws.Cells.Item(12, 11).Value = 1.044574317709982
$ws.Cells.Item(12, 12).Value = 1.051201782726125
$ws.Cells.Item(12, 13).Value = 1.056442221815254
$ws.Cells.Item(12, 14).Value = 1.016259360215746
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.030199000266889
$ws.Cells.Item(13, 4).Value = 1.041013500019671
$ws.Cells.Item(13, 5).Value = 1.047677891491108
$ws.Cells.Item(13, 6).Value = 1.052935653921408
$ws.Cells.Item(13, 9).Value = 1.037530551517119
$ws.Cells.Item(13, 10).Value = 1.036869263296747
$ws.Cells.Item(13, 11).Value = 1.044597255803492
$ws.Cells.Item(13, 12).Value = 1.051237119439711
$ws.Cells.Item(13, 13).Value = 1.056475773054478
$ws.Cells.Item(13, 14).Value = 1.016266574308085
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.030330930741868
$ws.Cells.Item(14, 4).Value = 1.041120487971787
$ws.Cells.Item(14, 5).Value = 1.0478254469517
$ws.Cells.Item(14, 6).Value = 1.053077144410716
$ws.Cells.Item(14, 9).Value = 1.037559133829048
$ws.Cells.Item(14, 10).Value = 1.036940095647148
$ws.Cells.Item(14, 11).Value = 1.044672412090825
$ws.Cells.Item(14, 12).Value = 1.051352930926477
$ws.Cells.Item(14, 13).Value = 1.056585725405713
$ws.Cells.Item(14, 14).Value = 1.016290210880199
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.030412234197924
$ws.Cells.Item(15, 4).Value = 1.041186422443791
$ws.Cells.Item(15, 5).Value = 1.047916398190913
$ws.Cells.Item(15, 6).Value = 1.053164352224735
$ws.Cells.Item(15, 9).Value = 1.037576723614904
$ws.Cells.Item(15, 10).Value = 1.036983737198985
$ws.Cells.Item(15, 11).Value = 1.044718716742917
$ws.Cells.Item(15, 12).Value = 1.0514243076977
$ws.Cells.Item(15, 13).Value = 1.056653485360712
$ws.Cells.Item(15, 14).Value = 1.01630477342681
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.030885709560016
$ws.Cells.Item(16, 4).Value = 1.041570427252564
$ws.Cells.Item(16, 5).Value = 1.048446343922372
$ws.Cells.Item(16, 6).Value = 1.053672408861561
$ws.Cells.Item(16, 9).Value = 1.037678786840953
$ws.Cells.Item(16, 10).Value = 1.037237739797533
$ws.Cells.Item(16, 11).Value = 1.044988204062897
$ws.Cells.Item(16, 12).Value = 1.051840078166211
$ws.Cells.Item(16, 13).Value = 1.057048100412779
$ws.Cells.Item(16, 14).Value = 1.01638952231223
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.031182930404172
$ws.Cells.Item(17, 4).Value = 1.04181151106331
$ws.Cells.Item(17, 5).Value = 1.048779265353023
$ws.Cells.Item(17, 6).Value = 1.053991511085384
$ws.Cells.Item(17, 9).Value = 1.037742528665645
$ws.Cells.Item(17, 10).Value = 1.03739705896827
$ws.Cells.Item(17, 11).Value = 1.045157222479425
$ws.Cells.Item(17, 12).Value = 1.052101166380423
$ws.Cells.Item(17, 13).Value = 1.057295826799296
$ws.Cells.Item(17, 14).Value = 1.016442672602994
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.031356373781407
$ws.Cells.Item(18, 4).Value = 1.041952205404617
$ws.Cells.Item(18, 5).Value = 1.048973632467291
$ws.Cells.Item(18, 6).Value = 1.054177785692993
$ws.Cells.Item(18, 9).Value = 1.037779607141204
$ws.Cells.Item(18, 10).Value = 1.037489983142192
$ws.Cells.Item(18, 11).Value = 1.045255799003902
$ws.Cells.Item(18, 12).Value = 1.052253557352082
$ws.Cells.Item(18, 13).Value = 1.057440391029643
$ws.Cells.Item(18, 14).Value = 1.016473670354007
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.031415526915464
$ws.Cells.Item(19, 4).Value = 1.04200019107046
$ws.Cells.Item(19, 5).Value = 1.049039937041534
$ws.Cells.Item(19, 6).Value = 1.05424132550331
$ws.Cells.Item(19, 9).Value = 1.037792232780449
$ws.Cells.Item(19, 10).Value = 1.037521667165645
$ws.Cells.Item(19, 11).Value = 1.04528940947135
$ws.Cells.Item(19, 12).Value = 1.052305536089298
$ws.Cells.Item(19, 13).Value = 1.057489695450311
$ws.Cells.Item(19, 14).Value = 1.016484239107929
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.031151033184994
$ws.Cells.Item(20, 4).Value = 1.041785637359516
$ws.Cells.Item(20, 5).Value = 1.048743527434067
$ws.Cells.Item(20, 6).Value = 1.053957259145909
$ws.Cells.Item(20, 9).Value = 1.037735700218088
$ws.Cells.Item(20, 10).Value = 1.037379965933706
$ws.Cells.Item(20, 11).Value = 1.045139089329541
$ws.Cells.Item(20, 12).Value = 1.052073143457184
$ws.Cells.Item(20, 13).Value = 1.05726924088834
$ws.Cells.Item(20, 14).Value = 1.016436970481326
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.030292077015305
$ws.Cells.Item(21, 4).Value = 1.041088979419845
$ws.Cells.Item(21, 5).Value = 1.047781987740713
$ws.Cells.Item(21, 6).Value = 1.053035472556694
$ws.Cells.Item(21, 9).Value = 1.037550721365681
$ws.Cells.Item(21, 10).Value = 1.036919237407439
$ws.Cells.Item(21, 11).Value = 1.044650280771714
$ws.Cells.Item(21, 12).Value = 1.051318822849647
$ws.Cells.Item(21, 13).Value = 1.056553344116808
$ws.Cells.Item(21, 14).Value = 1.016283250649154
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.029752746377899
$ws.Cells.Item(22, 4).Value = 1.040651645239306
$ws.Cells.Item(22, 5).Value = 1.047179062994523
$ws.Cells.Item(22, 6).Value = 1.052457254131403
$ws.Cells.Item(22, 9).Value = 1.037433512540227
$ws.Cells.Item(22, 10).Value = 1.036629531976608
$ws.Cells.Item(22, 11).Value = 1.044342875668797
$ws.Cells.Item(22, 12).Value = 1.050845488484686
$ws.Cells.Item(22, 13).Value = 1.056103871015322
$ws.Cells.Item(22, 14).Value = 1.016186568831511
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.030038585998031
$ws.Cells.Item(23, 4).Value = 1.040883419030678
$ws.Cells.Item(23, 5).Value = 1.047498529266436
$ws.Cells.Item(23, 6).Value = 1.052763650406487
$ws.Cells.Item(23, 9).Value = 1.037495732718661
$ws.Cells.Item(23, 10).Value = 1.036783112443936
$ws.Cells.Item(23, 11).Value = 1.044505843254962
$ws.Cells.Item(23, 12).Value = 1.051096322578789
$ws.Cells.Item(23, 13).Value = 1.056342083981145
$ws.Cells.Item(23, 14).Value = 1.016237824575642
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.031165445919987
$ws.Cells.Item(24, 4).Value = 1.041797328343437
$ws.Cells.Item(24, 5).Value = 1.048759675306439
$ws.Cells.Item(24, 6).Value = 1.053972735666987
$ws.Cells.Item(24, 9).Value = 1.037738786011944
$ws.Cells.Item(24, 10).Value = 1.037387689552804
$ws.Cells.Item(24, 11).Value = 1.045147282947619
$ws.Cells.Item(24, 12).Value = 1.05208580549308
$ws.Cells.Item(24, 13).Value = 1.057281253702248
$ws.Cells.Item(24, 14).Value = 1.016439547036819
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.032476215402256
$ws.Cells.Item(25, 4).Value = 1.042860774127882
$ws.Cells.Item(25, 5).Value = 1.050230205715813
$ws.Cells.Item(25, 6).Value = 1.055381597130682
$ws.Cells.Item(25, 9).Value = 1.038016869076986
$ws.Cells.Item(25, 10).Value = 1.038089105911908
$ws.Cells.Item(25, 11).Value = 1.045891277496112
$ws.Cells.Item(25, 12).Value = 1.053238064990846
$ws.Cells.Item(25, 13).Value = 1.058373832444395
$ws.Cells.Item(25, 14).Value = 1.016673479485006
